$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I0 / IF headers, matching H1's style (bold, bordered, centered) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: I = 1 (constant) and J = same value as H, except rows 24/25 ---
$values = @{
    2  = @(1, 5)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 6)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 5)
    15 = @(1, 7)
    16 = @(1, 6)
    17 = @(1, 4)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 3)
    22 = @(1, 5)
    23 = @(1, 6)
    24 = @(5, 7)
    25 = @(4, 6)
    26 = @(1, 2)
    27 = @(1, 1)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
